$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "nan" placeholder text in C8 so it becomes an empty inline string cell
$ws.Range("C8").Value = ""

# Add new row 9 with review data
$ws.Range("A9").Value = "parisk"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "Not too surprisingly"
$ws.Range("D9").Value = "CRT"
$ws.Range("E9").Value = "MET"
$ws.Range("F9").Value = "afe80f3f-3501-40b4-a3d0-1ad1f86c76ec"
$ws.Range("G9").Value = "r1BRfhiab_annotated.xlsx"
$ws.Range("H9").Value = "Not too surprisingly, the standard multiclass losses do not have the desired property, however approaches that reduce multi-class to binary classification at training time do, namely unnormalized models with penalized log Z (self-normalization), the NCE approach, as well as (the natural in the proposed setting) binary classification loss."
$ws.Range("I9").Value = "Correct"
